$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35: copy formatting pattern from row 29 (A=date s25, B/C/E/F s20, D s27, G s24)
$ws.Range("A29:G29").Copy()
$ws.Range("A35:G35").PasteSpecial(-4122)
$ws.Range("A35").Value = 43895
$ws.Range("B35").Value = "17:00-19:00 in class"
$ws.Range("C35").Value = "N.A."
$ws.Range("D35").Value = "Follow the lecture with professor"
$ws.Range("E35").Value = "Learnt more about Key Expert Practices and some testing methods"
$ws.Range("F35").Value = "Testing is not only important through developing process but also through reverse engineering process. Never tried to write any test cases before this quarter, but will be doing more from now on"
$ws.Range("G35").Value = "Good"

# Row 36: same formatting pattern as row 35 (copy from row 29 again)
$ws.Range("A29:G29").Copy()
$ws.Range("A36:G36").PasteSpecial(-4122)
$ws.Range("A36").Value = 43902
$ws.Range("B36").Value = "17:00-19:00 in class"
$ws.Range("C36").Value = "N.A."
$ws.Range("D36").Value = "Follow the lecture with professor"
$ws.Range("E36").Value = "Learnt 3 more Key Expert Practice and some advanced topics"
$ws.Range("F36").Value = "History of a project is an amazing topic which i didn’t treat as of valuable before. Git Blame is an interesting tool for developers. It won’t be that embarrassing and difficult now when i know how to use Git Blame properly through coding teamwork. jk. :) Also we should respect the wisdom in the old code and always leave code in a better place."
$ws.Range("G36").Value = "Average"

# Row 37: copy formatting pattern from row 33 (A=date s25, B/C/D/E/F s20, G s24)
$ws.Range("A33:G33").Copy()
$ws.Range("A37:G37").PasteSpecial(-4122)
$ws.Range("A37").Value = 43904
$ws.Range("B37").Value = "14:30-18:25"
$ws.Range("C37").Value = "Soobin, Marc"
$ws.Range("D37").Value = "Figure out issue for 2nd PR, assignment 6. Try to create new test cases for Glide."
$ws.Range("E37").Value = "Submit PR for issue proposal. Revise homework4 and start homework6"
$ws.Range("F37").Value = "We went through a hard time finding a proper issue to contribute. Since Glide doesn’t have many current (non-stale) issues and most of those open issues are either core-lib related or hard to reproduce on our laptop, it is kinda hard to choose one that we can fix before the ddl of our course. While we found one eventually, waiting for Kaj’s response and feeling nervous about the final"
$ws.Range("G37").Value = "Nervous about the final"
